$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of the shared string used by G5 ("If done, make plan for
# Monday" -> "...Tuesday") before it gets relocated below.
$ws.Range("G5").Value = "If done, make plan for Tuesday"

# "Failed attempts to make left-moving zombie": the right-hand mini table
# (originally in column G) got nudged one column over to H for most rows,
# while row 4 stayed put and row 6 just picked up the yellow highlight
# instead of actually moving.
$ws.Range("G2").Cut($ws.Range("H2"))
$ws.Range("G2").Clear()

$ws.Range("G3").Cut($ws.Range("H3"))
$ws.Range("G3").Clear()

$ws.Range("G5").Cut($ws.Range("H5"))
$ws.Range("G5").Clear()

# Row 6 cell never moved - it just got the same yellow fill used elsewhere
# in the sheet (style index 3: centered/wrapped + yellow fill).
$ws.Range("G6").Interior.Color = 65535

# Update the view: scrolled so B1 is the top-left visible cell, with the
# active selection now on H6 (following the column shift above).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H6").Select()
